$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# A2 is a text value that looks numeric ("  77"); force it to stay text
# (matches the source's inlineStr representation) without leaving a
# lingering number-format/style on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = " 77"
$ws.Range("A2").NumberFormat = "General"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = 213
$ws.Range("C2").Value = 1.01
$ws.Range("D2").Value = 178.6
$ws.Range("E2").Value = 0.94
$ws.Range("F2").Value = 118.7
$ws.Range("G2").Value = 0.73
$ws.Range("H2").Value = 0.66
$ws.Range("I2").Value = 0.39
$ws.Range("J2").Value = 0.08
$ws.Range("K2").Value = 0.08
$ws.Range("L2").Value = 140
$ws.Range("M2").Value = 149
$ws.Range("N2").Value = 83
$ws.Range("O2").Value = 16
$ws.Range("P2").Value = 17

# --- Row 3 ---
# A3 is likewise a text value with a leading space.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = " 15"
$ws.Range("A3").NumberFormat = "General"
$ws.Range("A3").Style = "Normal"
